$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.190.06'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '1.804.35'
$ws.Range('E3').Value = '  +0.63%  '
$ws.Range('D5').Value = "'223.57"
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = "'0.554"
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D8').Value = "'32.40"
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('D9').Value = "'0.287"
$ws.Range('E9').Value = '  +2.34%  '
$ws.Range('D10').Value = "'0.0723"
$ws.Range('E10').Value = '  +4.92%  '
$ws.Range('D11').Value = "'0.0927"
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('D12').Value = '2.063.16'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').Value = "'11.02"
$ws.Range('E13').Value = '  +1.48%  '
$ws.Range('D14').Value = '1.793.83'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = "'0.632"
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').Value = '34.209.25'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = "'248.18"
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').Value = '0.0₃0791'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').Value = "'10.97"
$ws.Range('E21').Value = '  +5.04%  '
$ws.Range('D22').Value = "'1.00"
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('D25').Value = "'159.55"
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').Value = "'16.67"
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('E27').Value = '  +0.65%  '
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('E32').Value = '  +1.92%  '
$ws.Range('D33').Value = "'3.53"
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').Value = '1.418.65'
$ws.Range('E35').Value = '  -1.30%  '
$ws.Range('D36').Value = "'0.655"
$ws.Range('E36').Value = '  +2.91%  '
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('D39').Value = "'0.947"
$ws.Range('E39').Value = '  +3.85%  '
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('D41').Value = "'2.36"
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').Value = "'2.73"
$ws.Range('E42').Value = '  -2.62%  '
$ws.Range('D43').Value = "'2.16"
$ws.Range('E43').Value = '  +3.50%  '
$ws.Range('D44').Value = "'5.95"
$ws.Range('E44').Value = '  -1.38%  '
$ws.Range('D45').Value = "'108.06"
$ws.Range('E45').Value = '  +3.91%  '
$ws.Range('D46').Value = "'0.0497"
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('D47').Value = '1.961.20'
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('D48').Value = "'1.05"
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('D49').Value = "'12.04"
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('E51').Value = '  +3.92%  '
